$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range("A21").Value = "R5EJWA"
$ws.Range("B21").Value = "2024-12-25 19:39:43"
$ws.Range("C21").Value = "GET /"
$ws.Range("D21").Value = 401
$ws.Range("E21").Value = $False
$ws.Range("F21").Value = "Eitss... mau ngapain? Akses terbatas!"

# Row 22
$ws.Range("A22").Value = "C8HEHX"
$ws.Range("B22").Value = "2024-12-25 19:40:04"
$ws.Range("C22").Value = "GET /checkmodel"
$ws.Range("D22").Value = 200
$ws.Range("E22").Value = $True
$ws.Range("F22").Value = "OK`n###`nLog Config LLM:{'last_update': '2024-12-25 19:29:28', 'llm': 'openai', 'model_llm': 'gpt-4o', 'embbeder': 'openai', 'model_embedder': 'text-embedding-3-large', 'chunk_size': 900, 'chunk_overlap': 100, 'total_chunks': 177}"

# Row 23
$ws.Range("A23").Value = "J4CNYW"
$ws.Range("B23").Value = "2024-12-25 19:40:36"
$ws.Range("C23").Value = "POST /setup"
$ws.Range("D23").Value = 200
$ws.Range("E23").Value = $True
$ws.Range("F23").Value = "Proses penyiapan dokumen berhasil diselesaikan dan embeddings berhasil disimpan pada vector database.`n###`nllm:openai`n###`nmodel_llm:gpt-4o`n###`nembbeder:openai`n###`nmodel_embedder:text-embedding-3-large`n###`nchunk_size:1000`n###`nchunk_overlap:200`n###`ntotal_chunks:173"

# Row 24
$ws.Range("A24").Value = "GIJVMN"
$ws.Range("B24").Value = "2024-12-25 19:49:44"
$ws.Range("C24").Value = "POST /setup"
$ws.Range("D24").Value = 200
$ws.Range("E24").Value = $True
$ws.Range("F24").Value = "Proses penyiapan dokumen berhasil diselesaikan dan embeddings berhasil disimpan pada vector database.`n###`nllm:openai`n###`nmodel_llm:gpt-4o`n###`nembbeder:openai`n###`nmodel_embedder:text-embedding-3-large`n###`nchunk_size:1000`n###`nchunk_overlap:200`n###`ntotal_chunks:173"

# Row 25
$ws.Range("A25").Value = "UH625V"
$ws.Range("B25").Value = "2024-12-25 19:53:00"
$ws.Range("C25").Value = "POST /setup"
$ws.Range("D25").Value = 400
$ws.Range("E25").Value = $False
$ws.Range("F25").Value = "Model Embedder untuk 'openai' harus salah satu dari ['text-embedding-3-large', 'text-embedding-3-small']."

# Row 26
$ws.Range("A26").Value = "45988Q"
$ws.Range("B26").Value = "2024-12-25 19:53:07"
$ws.Range("C26").Value = "POST /setup"
$ws.Range("D26").Value = 200
$ws.Range("E26").Value = $True
$ws.Range("F26").Value = "Proses penyiapan dokumen berhasil diselesaikan dan embeddings berhasil disimpan pada vector database.`n###`nllm:openai`n###`nmodel_llm:gpt-4o`n###`nembbeder:openai`n###`nmodel_embedder:text-embedding-3-large`n###`nchunk_size:1000`n###`nchunk_overlap:200`n###`ntotal_chunks:173"
